# Generate Report for Handoff
# ----------------------------
# The localization source file was renamed from
#   a5c55c48-9368-488a-9ab6-29f00c1b2b14.md
# to
#   08e7e948-2bd4-42b7-87f1-f490872ad413.md
# and a fresh handoff pass produced new handoff-file names / timestamps
# for both the zh-cn and de-de targets. Update every cell (and the
# matching hyperlink display text) that references the old id/time so
# the workbook reflects the new handoff report.

$wb = $excel.ActiveWorkbook

$oldId = "a5c55c48-9368-488a-9ab6-29f00c1b2b14"
$newId = "08e7e948-2bd4-42b7-87f1-f490872ad413"
$oldHash = "a7f8746f57ef2c44ea2ed3336e92a497d619cf90"
$newHash = "829ca35757c53b82192e252ea8bc2e0ce3379203"

$newMdName = "$newId.md"

$newZhHandoffFile = "$newId.$newHash.zh-cn.xlf"
$newZhHandoffTime = "2016-02-22 14:14:48"

$newDeHandoffFile = "$newId.$newHash.de-de.xlf"
$newDeHandoffTime = "2016-02-22 14:15:04"

function Set-CellAndHyperlink {
    param($Worksheet, $CellRef, $NewText)

    $range = $Worksheet.Range($CellRef)
    $range.Value = $NewText

    # Keep the hyperlink's visible text (TextToDisplay / display attribute)
    # in sync with the new cell text, for every hyperlink anchored here.
    for ($i = 1; $i -le $Worksheet.Hyperlinks.Count; $i++) {
        $hl = $Worksheet.Hyperlinks.Item($i)
        if ($hl.Range.Address() -eq $range.Address()) {
            $hl.TextToDisplay = $NewText
        }
    }
}

# --- Overview sheet: A2 holds the source file name -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $wsOverview "A2" $newMdName

# --- zh-cn sheet: A2 source file name, C2 handoff file, D2 handoff time --
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $wsZh "A2" $newMdName
Set-CellAndHyperlink $wsZh "C2" $newZhHandoffFile
$wsZh.Range("D2").Value = $newZhHandoffTime

# --- de-de sheet: A2 source file name, C2 handoff file, D2 handoff time --
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $wsDe "A2" $newMdName
Set-CellAndHyperlink $wsDe "C2" $newDeHandoffFile
$wsDe.Range("D2").Value = $newDeHandoffTime
